# Edit script for "Hortaliza, Vega Monumental Concepción - Betarraga.xlsx"
#
# The underlying data block (rows 73..265, columns D..R) is shifted DOWN by
# two rows: sheet row (r) ends up holding what used to live in row (r-2),
# for r = 75..265. The two rows that "fall off" the end of the original
# block (old rows 264 and 265) become two brand-new rows, 266 and 267
# (with the same Mercado ID / Mercado / Región values repeated in A:C, as
# is the pattern for every row in this sheet). Rows 73 and 74 keep all of
# their original data except the date in column D, which becomes a new
# value (44607) not derived from any existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow   = 73
$endRow     = 265
$numOldRows = $endRow - $startRow + 1      # 193
$numNewRows = $numOldRows + 2              # 195  (rows 73..267)
$numCols    = 15                            # columns D..R

# Snapshot the original D:R block for rows 73..265 in one COM round-trip.
$oldRange = $ws.Range("D" + $startRow + ":R" + $endRow)
$old = $oldRange.Value2

# Build the replacement block (0-based .NET array; row 0 -> sheet row 73).
$new = New-Object 'object[,]' $numNewRows, $numCols

# Rows 73 and 74 (new-block rows 0,1) keep their own original content.
for ($c = 0; $c -lt $numCols; $c++) {
    $new[0, $c] = $old[1, $c + 1]
    $new[1, $c] = $old[2, $c + 1]
}
# ...except column D (block column 0), which takes a brand-new date value.
$new[0, 0] = 44607
$new[1, 0] = 44607

# Rows 75..267 (new-block rows 2..194) take the content that used to sit
# two rows above them (old rows 73..265), i.e. new-block row i = old row
# (i-1) in 1-based terms.
for ($i = 2; $i -lt $numNewRows; $i++) {
    $oldRowIdx = $i - 1
    for ($c = 0; $c -lt $numCols; $c++) {
        $new[$i, $c] = $old[$oldRowIdx, $c + 1]
    }
}

$newRange = $ws.Range("D" + $startRow + ":R" + ($startRow + $numNewRows - 1))
$newRange.Value2 = $new

# The two brand-new rows (266, 267) still need A:C (Mercado ID / Mercado /
# Región) filled in -- every row in the sheet repeats the same three
# values, so copy them down from an existing row.
$ws.Range("A266:A267").Value2 = $ws.Range("A265").Value2
$ws.Range("B266:B267").Value2 = $ws.Range("B265").Value2
$ws.Range("C266:C267").Value2 = $ws.Range("C265").Value2

# Column D is a date column formatted specially (numFmtId 165); new rows
# need that same number format applied (existing rows already have it).
$dateFormat = $ws.Range("D" + $startRow).NumberFormat
$ws.Range("D266:D267").NumberFormat = $dateFormat
